# AH_cfDNA_v3_VAR.docx - Release 1.0 template changes
#
# 1) Assay version reference updated from "v1" to "v3" in the methods
#    paragraph, and the sentence is split across two runs at that point.
# 2) The cached "Reported" date field result text is bumped from
#    1-Nov-2023 to 15-Nov-2023.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: " v1) and sequenced on an Illumina NextSeq500 ..." -> " v3)" +
#         ") and sequenced on an Illumina NextSeq500 ..." (two runs).
# ---------------------------------------------------------------------

# Locate the exact text we need to touch.
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute(" v1) and sequenced on an Illumina NextSeq500 with 150 bp paired end reads. ",
                              $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the assay-version sentence to edit."
}
$sentenceStart = $rng1.Start
$sentenceEnd = $rng1.End

# The run boundary that must be preserved immediately after this sentence
# (the following run, e.g. "A customised CLC ...", must stay untouched).
# Drop a throwaway bookmark there first so the text mutation below can't
# coalesce across it.
$boundary = $d.Range($sentenceEnd, $sentenceEnd)
$d.Bookmarks.Add("zzEditBoundary1", $boundary) | Out-Null

# Replace just the "1" in "v1" with "3".
$numStart = $sentenceStart + 2
$numEnd = $numStart + 1
$rngNum = $d.Range($numStart, $numEnd)
if ($rngNum.Text -ne "1") {
    throw "Unexpected character at version-number position: [$($rngNum.Text)]"
}
$rngNum.Text = "3"

# Re-touch the remainder of the sentence (a no-op formatting round trip)
# so it becomes its own run, split off from " v3".
$rngRemainder = $d.Range($numEnd, $sentenceEnd)
$rngRemainder.Font.Bold = 1
$rngRemainder.Font.Bold = 0

# Clean up the helper bookmark.
$d.Bookmarks("zzEditBoundary1").Delete()

# ---------------------------------------------------------------------
# Edit 2: cached "Reported" date field text, 1-Nov-2023 -> 15-Nov-2023.
# ---------------------------------------------------------------------

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("1-Nov-2023", $true, $false, $false, $false, $false,
                    $true, 1, $false, "15-Nov-2023", 2) | Out-Null
